$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "খাতা/পত্রের সংখ্যা" (quantity) counts that were left blank,
# which in turn makes the dependent formulas in column I (and the I32 total)
# recalculate to their real amounts instead of 0.
$ws.Range("G16").Value = 27
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1
